# Auto-generated edit script: updates cryptocurrency price/volume data
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.623.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.63%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.683.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -7.47%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.58%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.64%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.670.57"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -7.60%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -10.35%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.06%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.700"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -13.55%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -9.29%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000291"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -14.00%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -11.62%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.269.56"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.47%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.696.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.06%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -10.20%  "
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.48%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -10.37%  "
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -10.91%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.558.19"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.55%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -11.58%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.47"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.01%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -10.31%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -10.76%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -11.35%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.46%  "
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.54%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -12.57%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -12.10%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.57"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -10.51%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.47"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.08%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -11.78%  "
# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -10.41%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "64.89"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.72%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -12.57%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "597.24"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.81%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0883"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -14.90%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.02%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.395"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.35%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.11%  "
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -10.13%  "
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -13.12%  "
# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.05%  "
# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.92"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -11.77%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0435"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -10.77%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.16"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -13.48%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.791.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.80%  "
# Row 49
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.37%  "
# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.133"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -11.45%  "
# Row 51
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.87%  "

Write-Output "Applied cryptos list update: $($wb.ActiveSheet.Name)"
